# Generate Report for Archive
#
# 1. "Ready for handoff" -> "In Translation" (shared-string text used by the
#    Status columns on all three sheets).
# 2. The Status columns' widths shrink to match the new (shorter) text
#    (17.2159881591797 -> 13.4101845877511 raw OOXML width units, i.e. the
#    closest value the Excel column-width pixel grid can represent).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update the status text on every sheet that references it ---------
# Overview sheet: zh-cn status is column E, de-de status is column F.
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# zh-cn / de-de detail sheets: Status is column C.
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Shrink the Status columns to fit the new, shorter text ------------
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
